$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-blank row 48 with its log entry
$ws.Range("A48").Value = 43080
$ws.Range("B48").Value = 0.45833333333333331
$ws.Range("C48").Value = 0.48958333333333331

# Move the active selection to C49 (was F43)
[void]$ws.Range("C49").Select()
